# Edit script generated for Config.xlsx business rules update
$wb = $excel.ActiveWorkbook

# --- Settings sheet: add new Asset abbreviation rows (33-35) ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A33").Value = 'CSC'
$wsSettings.Range("B33").Value = 'Child Support Calcs'
$wsSettings.Range("A34").Value = 'CSPH'
$wsSettings.Range("B34").Value = 'Child Support Payment History'
$wsSettings.Range("A35").Value = 'SSBL'
$wsSettings.Range("B35").Value = 'Social Security Benefits Letter'

# --- Findings sheet: add new business rule rows ---
$wsFindings = $wb.Worksheets.Item("Findings")
$wsFindings.Range("A169").Value = 'CSPH_ParentNameMatch'
$wsFindings.Range("B169").Value = 'Receiving Parent listed on the Child Support Payment History, Child Support Calcs, and ICW do not match.'
$wsFindings.Range("A170").Value = 'CSPH_FrequencyCheck'
$wsFindings.Range("B170").Value = 'Unable to determine frequency of payments from the Child Support Payment History.'
$wsFindings.Range("A171").Value = 'CSPH_DateCheck'
$wsFindings.Range("B171").Value = 'Date on the Child Support Payment History is over 120 days old. '
$wsFindings.Range("A172").Value = 'CSPH_ChildrenNamesCheck'
$wsFindings.Range("B172").Value = 'Child''s Name listed on the Child Support Payment History and Child Support Calcs do not match.'
$wsFindings.Range("A173").Value = 'CSPH_PaymentAmountCheck'
$wsFindings.Range("B173").Value = 'Payment Amount listed on the Child Support Payment History and Child Support Calcs do not match.'
$wsFindings.Range("A174").Value = 'CSPH_PaymentsInYearCheck'
$wsFindings.Range("B174").Value = 'Payments in a Year listed on the Child Support Payment History and Child Support Calcs do not match.'
$wsFindings.Range("A175").Value = 'CSPH_TotalPerYearCheck'
$wsFindings.Range("B175").Value = 'Total Child Support Payment in a Year listed on Child Support Calcs and the ICW do not match.'
$wsFindings.Range("A176").Value = 'CSPH_StartDateCheck'
$wsFindings.Range("B176").Value = 'Start Date of Payments listed on the Child Support Payment History and Child Support Calcs do not match.'
$wsFindings.Range("A177").Value = 'CSPH_EndDateCheck'
$wsFindings.Range("B177").Value = 'End Date of Payments listed on the Child Support Payment History and Child Support Calcs do not match.'
$wsFindings.Range("A178").Value = 'CSPH_SumPaymentsCheck'
$wsFindings.Range("B178").Value = 'Sum of Payments calculated from Child Support Payment History and listed on Child Support Calcs do not match.'
$wsFindings.Range("A179").Value = 'CSPH_PerYearPaymentsInICW'
$wsFindings.Range("B179").Value = 'Total Child Support Payment in a Year listed on Child Support Calcs and the ICW do not match.'
$wsFindings.Range("A181").Value = 'SSBL_EmployeeCheck'
$wsFindings.Range("B181").Value = 'Applicant Name listed on the Social Security Benefits Letter and ICW do not match. '
$wsFindings.Range("A182").Value = 'SSBL_DateCheck'
$wsFindings.Range("B182").Value = 'Date of the Social Security Benefits Letter is over 120 days in the past.'
$wsFindings.Range("A183").Value = 'SSBL_BNCNumberCheck'
$wsFindings.Range("B183").Value = 'BNC# listed on the Social Security Benefits Letter does not match on every page.'
$wsFindings.Range("A184").Value = 'SSBL_MonthlyAmountCheck'
$wsFindings.Range("B184").Value = 'Monthly Amount listed on the Social Security Benefits Letter and the ICW do not match.'
$wsFindings.Range("A185").Value = 'SSBL_YearlyAmountCheck'
$wsFindings.Range("B185").Value = 'Yearly Amount calculated from the Social Security Benefits Letter does not math the ICW.'
$wsFindings.Range("A186").Value = 'SSBL_COLALetterCheck'
$wsFindings.Range("B186").Value = 'COLA Letter is missing from the application.'
$wsFindings.Range("A187").Value = 'SSBL_PressReleaseYearCheck'
$wsFindings.Range("B187").Value = 'Press Release Year from the COLA Letter is not next year as expected.'
$wsFindings.Range("A188").Value = 'SSBL_ApplicantNameCheckCola'
$wsFindings.Range("B188").Value = 'Applicant Name listed on the Social Security Benefits Letter and COLA Calculator do not match. '
$wsFindings.Range("A189").Value = 'SSBL_MonthlyAmountCheckCola'
$wsFindings.Range("B189").Value = 'Monthly Amount listed on the Social Security Benefits Letter and the COLA Calculator do not match.'
$wsFindings.Range("A190").Value = 'SSBL_MIMonthCheck'
$wsFindings.Range("B190").Value = 'Move In Month listed on the COLA Calculator does not match that as the next full months following the Move In Date on the Application Summary.'
$wsFindings.Range("A191").Value = 'SSBL_TotalPerYearMatchICW'
$wsFindings.Range("B191").Value = 'Total per year listed on the COLA Calculator and ICW do not match.'

# --- Restore selections to match authored state ---
$null = $wsSettings.Activate()
$null = $wsSettings.Range("A38").Select()

$null = $wsFindings.Activate()
$null = $wsFindings.Range("B186").Select()
